# AssureTestData.xlsx - "LoginTest" sheet edit
#
# Before:
#   Row1: (blank)
#   Row2: A2=LoginTest
#   Row3: A3=InputUsername   B3=InputPassword  C3=ExpectedUsername   (bold/yellow header)
#   Row4: A4=<email>         B4=<password>     C4=<greeting>
#   Row5: A5=<email>         B5=<password>     C5=<greeting>
#   Row6: A6=<email>         B6=<password>     C6=<greeting>
#
# After:
#   Row1: A1=LoginTest
#   Row2: A2=InputUsername   B2=InputPassword  C2=ExpectedUsername   (bold/yellow header)
#   Row3: A3=<email>         B3=<password>     C3=<greeting>
#   Row4: A4=<email>         B4=<password>     C4=<greeting>
#   Row5: A5=<email>         B5=N0rthg4t311    C5=<greeting>
#
# i.e. the leading blank row is removed (shifting everything up by one),
# and the last data row gets a new, distinct password value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Deleting the empty first row shifts rows 2-6 up to 1-5, carrying their
# existing formatting (the bold/yellow header row) along with them.
$ws.Rows("1:1").Delete() | Out-Null

# The final row's password is now a distinct value rather than a repeat
# of the earlier rows.
$ws.Range("B5").Value = "N0rthg4t311"

# Land the cursor back on the top-left cell.
$ws.Range("A1").Select()
